$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '28.686.56'
$ws.Cells.Item(2, 5).Value = '  +2.48%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.905.06'
$ws.Cells.Item(3, 5).Value = '  +2.49%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +2.70%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '320.18'

# Row 6
$ws.Cells.Item(6, 5).Value = '  +2.61%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.5204'
$ws.Cells.Item(7, 5).Value = '  +1.15%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3954'
$ws.Cells.Item(8, 5).Value = '  +3.13%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.08347'
$ws.Cells.Item(9, 5).Value = '  +1.47%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +2.16%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +2.42%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '6.301'
$ws.Cells.Item(12, 5).Value = '  +1.71%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '1.913.06'
$ws.Cells.Item(13, 5).Value = '  +2.70%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '20.68'
$ws.Cells.Item(14, 5).Value = '  +0.61%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '7.315'
$ws.Cells.Item(15, 5).Value = '  +0.83%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '1.032'
$ws.Cells.Item(16, 5).Value = '  +2.86%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.00001115'
$ws.Cells.Item(17, 5).Value = '  +1.67%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '91.61'
$ws.Cells.Item(18, 5).Value = '  +1.11%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.06799'
$ws.Cells.Item(19, 5).Value = '  +2.31%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '18.02'
$ws.Cells.Item(20, 5).Value = '  +1.88%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '1.029'
$ws.Cells.Item(21, 5).Value = '  +2.60%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.092'
$ws.Cells.Item(22, 5).Value = '  +1.40%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '28.735.39'
$ws.Cells.Item(23, 5).Value = '  +2.56%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '11.26'
$ws.Cells.Item(24, 5).Value = '  +1.66%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.294'
$ws.Cells.Item(25, 5).Value = '  +1.31%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '2.129.60'
$ws.Cells.Item(26, 5).Value = '  +2.73%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '162.68'
$ws.Cells.Item(27, 5).Value = '  +3.06%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '21.00'
$ws.Cells.Item(28, 5).Value = '  +2.70%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.459'
$ws.Cells.Item(29, 5).Value = '  -2.01%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '127.54'
$ws.Cells.Item(30, 5).Value = '  +2.30%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.1070'
$ws.Cells.Item(31, 5).Value = '  +0.72%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.051'
$ws.Cells.Item(32, 5).Value = '  +2.05%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '6.001'
$ws.Cells.Item(33, 5).Value = '  +0.93%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '3.684'
$ws.Cells.Item(34, 5).Value = '  +2.41%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'FraxShare'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '9.451'
$ws.Cells.Item(35, 5).Value = '  +0.87%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'VeChain'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.02470'
$ws.Cells.Item(36, 5).Value = '  +2.24%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.06639'
$ws.Cells.Item(37, 5).Value = '  +2.12%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.2222'
$ws.Cells.Item(38, 5).Value = '  +2.09%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.6563'
$ws.Cells.Item(39, 5).Value = '  -0.09%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.262'
$ws.Cells.Item(40, 5).Value = '  +3.68%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.196'
$ws.Cells.Item(41, 5).Value = '  +0.10%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '5.026'
$ws.Cells.Item(42, 5).Value = '  +0.75%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '11.14'
$ws.Cells.Item(43, 5).Value = '  -0.20%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.6179'
$ws.Cells.Item(44, 5).Value = '  +0.23%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '13.20'
$ws.Cells.Item(45, 5).Value = '  +1.41%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '3.761'
$ws.Cells.Item(46, 5).Value = '  +2.51%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.303'
$ws.Cells.Item(47, 5).Value = '  +1.62%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.018'
$ws.Cells.Item(48, 5).Value = '  +0.68%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +2.00%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '122.58'
$ws.Cells.Item(50, 5).Value = '  +1.61%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.06960'
$ws.Cells.Item(51, 5).Value = '  +2.54%  '
